# Fix checklist answer validation error
#
# The "action" data-validation list was anchored to column C (and G) on the
# "Checklist" sheet. A new "Section" column needs to sit right after "S/N"
# (i.e. become the new column B), pushing "Action" and everything after it
# one column to the right. Excel automatically slides the dataValidation
# sqref ranges (C/G -> D/H, F -> G) and defined-name references along with
# a real column insert, which is exactly what was broken/fixed upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert a blank column before the current column B ("Action"), shifting
# Action..Due Date from B:J to C:K.
$ws.Columns.Item(2).Insert()

# Label the newly inserted column.
$ws.Cells.Item(1, 2).Value = "Section"

# Leave the same active cell selection behind as the authored edit.
$ws.Range("C4").Select() | Out-Null
